# Update the B15 result value with the newly recorded experiment result.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B15").Value = 0.9212

# Move/save the active selection to C15 (matches the author's last-selected cell).
$ws.Range("C15").Select()
